$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the touched cells to remain plain text (the crypto price/volume
# columns store values like "0.9990" or "30.169.01" as text, and a plain
# COM .Value assignment would otherwise coerce them into numbers/dates).
$touched = $ws.Range("D2:E51")
$touched.NumberFormat = "@"

$ws.Range("D2").Value = '30.169.01'
$ws.Range("E2").Value = '  +4.42%  '
$ws.Range("D3").Value = '1.905.76'
$ws.Range("E3").Value = '  +5.09%  '
$ws.Range("D4").Value = '0.9990'
$ws.Range("E4").Value = '  -0.13%  '
$ws.Range("D5").Value = '252.31'
$ws.Range("E5").Value = '  +2.39%  '
$ws.Range("E6").Value = '  -0.15%  '
$ws.Range("D7").Value = '0.5087'
$ws.Range("E7").Value = '  +3.43%  '
$ws.Range("D8").Value = '45.11'
$ws.Range("E8").Value = '  +2.80%  '
$ws.Range("D9").Value = '0.3021'
$ws.Range("E9").Value = '  +8.92%  '
$ws.Range("D10").Value = '0.06806'
$ws.Range("E10").Value = '  +6.46%  '
$ws.Range("D11").Value = '1.905.21'
$ws.Range("E11").Value = '  +5.11%  '
$ws.Range("D12").Value = '17.34'
$ws.Range("E12").Value = '  +3.64%  '
$ws.Range("D13").Value = '0.07326'
$ws.Range("E13").Value = '  +3.60%  '
$ws.Range("D14").Value = '0.6926'
$ws.Range("E14").Value = '  +7.97%  '
$ws.Range("D15").Value = '87.04'
$ws.Range("E15").Value = '  +3.83%  '
$ws.Range("D16").Value = '4.919'
$ws.Range("E16").Value = '  +5.38%  '
$ws.Range("D17").Value = '30.148.15'
$ws.Range("E17").Value = '  +4.26%  '
$ws.Range("D18").Value = '0.000008239'
$ws.Range("E18").Value = '  +12.73%  '
$ws.Range("D19").Value = '0.9987'
$ws.Range("E19").Value = '  -0.07%  '
$ws.Range("D20").Value = '13.06'
$ws.Range("E20").Value = '  +6.82%  '
$ws.Range("D21").Value = '2.150.67'
$ws.Range("E21").Value = '  +5.33%  '
$ws.Range("D22").Value = '0.9985'
$ws.Range("E22").Value = '  -0.25%  '
$ws.Range("D23").Value = '4.829'
$ws.Range("E23").Value = '  +5.78%  '
$ws.Range("D24").Value = '5.747'
$ws.Range("E24").Value = '  +7.63%  '
$ws.Range("D25").Value = '9.372'
$ws.Range("E25").Value = '  +6.59%  '
$ws.Range("D26").Value = '148.38'
$ws.Range("E26").Value = '  +3.44%  '
$ws.Range("D27").Value = '134.44'
$ws.Range("E27").Value = '  +4.40%  '
$ws.Range("D28").Value = '17.15'
$ws.Range("E28").Value = '  +4.38%  '
$ws.Range("D29").Value = '2.007'
$ws.Range("E29").Value = '  +6.66%  '
$ws.Range("E30").Value = '  -0.57%  '
$ws.Range("D31").Value = '4.296'
$ws.Range("E31").Value = '  +4.28%  '
$ws.Range("D32").Value = '0.08900'
$ws.Range("D33").Value = '4.007'
$ws.Range("E33").Value = '  +6.22%  '
$ws.Range("D34").Value = '0.05203'
$ws.Range("E34").Value = '  +5.85%  '
$ws.Range("D35").Value = '1.149'
$ws.Range("E35").Value = '  +5.07%  '
$ws.Range("D36").Value = '0.7207'
$ws.Range("E36").Value = '  +7.63%  '
$ws.Range("D37").Value = '2.685'
$ws.Range("E37").Value = '  -0.41%  '
$ws.Range("E38").Value = '  +2.66%  '
$ws.Range("D39").Value = '2.293'
$ws.Range("E39").Value = '  +0.89%  '
$ws.Range("D40").Value = '0.9645'
$ws.Range("E40").Value = '  +1.59%  '
$ws.Range("D41").Value = '0.01692'
$ws.Range("E41").Value = '  +6.83%  '
$ws.Range("D42").Value = '6.096'
$ws.Range("E42").Value = '  -0.72%  '
$ws.Range("D43").Value = '0.4325'
$ws.Range("D44").Value = '105.29'
$ws.Range("E44").Value = '  +5.01%  '
$ws.Range("D45").Value = '0.9989'
$ws.Range("E45").Value = '  -0.15%  '
$ws.Range("D46").Value = '7.686'
$ws.Range("E46").Value = '  +8.09%  '
$ws.Range("D47").Value = '0.1282'
$ws.Range("E47").Value = '  +5.12%  '
$ws.Range("E48").Value = '  +4.35%  '
$ws.Range("D49").Value = '33.50'
$ws.Range("E49").Value = '  +6.03%  '
$ws.Range("D50").Value = '8.421'
$ws.Range("E50").Value = '  +3.51%  '
$ws.Range("D51").Value = '0.3824'
$ws.Range("E51").Value = '  +5.96%  '

# Restore the default (unstyled) cell format so only the values change,
# matching the original workbook formatting.
$touched.Style = "Normal"

